$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3494.7297
$ws.Range("J17").Value = 3735.4412
$ws.Range("L17").Value = 11206.3236
$ws.Range("N17").Value = -11542.3236
$ws.Range("H28").Value = 56724.89
$ws.Range("I28").Value = 91389.37
$ws.Range("K28").Value = 91389.37
$ws.Range("M28").Value = -90904.37
$ws.Range("H40").Value = 3408.5833
$ws.Range("I40").Value = 1890.7142
$ws.Range("K40").Value = 1890.7142
$ws.Range("M40").Value = -1715.7142
$ws.Range("H62").Value = 5440.1333
$ws.Range("I62").Value = 4717.3687
$ws.Range("J62").Value = 6688.5454
$ws.Range("K62").Value = 4717.3687
$ws.Range("L62").Value = 6688.5454
$ws.Range("M62").Value = -4093.3687
$ws.Range("N62").Value = -7936.5454
$ws.Range("H64").Value = 8427.571
$ws.Range("J64").Value = 8500.25
$ws.Range("L64").Value = 8500.25
$ws.Range("N64").Value = -8996.25
$ws.Range("H65").Value = 5440.1333
$ws.Range("I65").Value = 4717.3687
$ws.Range("J65").Value = 6688.5454
$ws.Range("K65").Value = 23586.8435
$ws.Range("L65").Value = 33442.727
$ws.Range("M65").Value = -20466.8435
$ws.Range("N65").Value = -39682.727
$ws.Range("H67").Value = 8427.571
$ws.Range("J67").Value = 8500.25
$ws.Range("L67").Value = 8500.25
$ws.Range("N67").Value = -10216.25
$ws.Range("H80").Value = 2207
$ws.Range("I80").Value = 550.55554
$ws.Range("J80").Value = 3353.7693
$ws.Range("K80").Value = 1651.66662
$ws.Range("L80").Value = 10061.3079
$ws.Range("M80").Value = -653.66662
$ws.Range("N80").Value = -12057.3079
$ws.Range("H83").Value = 2207
$ws.Range("I83").Value = 550.55554
$ws.Range("J83").Value = 3353.7693
$ws.Range("K83").Value = 4954.99986
$ws.Range("L83").Value = 30183.9237
$ws.Range("M83").Value = 37.0001400000001
$ws.Range("N83").Value = -40167.9237
$ws.Range("H86").Value = 5387.273
$ws.Range("J86").Value = 6479.8
$ws.Range("L86").Value = 6479.8
$ws.Range("N86").Value = -8725.799999999999
$ws.Range("H89").Value = 5387.273
$ws.Range("J89").Value = 6479.8
$ws.Range("L89").Value = 32399
$ws.Range("N89").Value = -43631
$ws.Range("H106").Value = 8865.333000000001
$ws.Range("I106").Value = 3735.8948
$ws.Range("K106").Value = 3735.8948
$ws.Range("M106").Value = -3104.8948
$ws.Range("H107").Value = 4029.15
$ws.Range("I107").Value = 4443.2144
$ws.Range("J107").Value = 3063
$ws.Range("K107").Value = 4443.2144
$ws.Range("L107").Value = 3063
$ws.Range("M107").Value = -2523.2144
$ws.Range("N107").Value = -6903
$ws.Range("H113").Value = 6231.1113
$ws.Range("I113").Value = 3495
$ws.Range("J113").Value = 7012.857
$ws.Range("K113").Value = 3495
$ws.Range("L113").Value = 7012.857
$ws.Range("M113").Value = -241
$ws.Range("N113").Value = -13520.857
$ws.Range("H135").Value = 1859.579
$ws.Range("I135").Value = 843
$ws.Range("J135").Value = 10500.5
$ws.Range("K135").Value = 7587
$ws.Range("L135").Value = 94504.5
$ws.Range("M135").Value = -5052
$ws.Range("N135").Value = -99574.5
$ws.Range("H138").Value = 4108.227
$ws.Range("I138").Value = 7965.3335
$ws.Range("K138").Value = 23896.0005
$ws.Range("M138").Value = -18756.0005
$ws.Range("H141").Value = 1049.3
$ws.Range("I141").Value = 1055.3334
$ws.Range("K141").Value = 3166.0002
$ws.Range("M141").Value = 2013.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4756.4546
$ws.Range("I45").Value = 1629.7142
$ws.Range("J45").Value = 10228.25
$ws.Range("K45").Value = 1629.7142
$ws.Range("L45").Value = 10228.25
$ws.Range("M45").Value = -1252.7142
$ws.Range("N45").Value = -10982.25
$ws.Range("H55").Value = 4413.5
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H97").Value = 1035.8334
$ws.Range("I97").Value = 1671.6666
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 1671.6666
$ws.Range("L97").Value = 400
$ws.Range("M97").Value = -1175.6666
$ws.Range("N97").Value = -1392
$ws.Range("H132").Value = 5046.8
$ws.Range("I132").Value = 3047
$ws.Range("K132").Value = 9141
$ws.Range("M132").Value = -6611

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 17788.7
$ws.Range("I105").Value = 16805.357
$ws.Range("K105").Value = 16805.357
$ws.Range("M105").Value = -15058.357
$ws.Range("H107").Value = 1264.2858
$ws.Range("I107").Value = 972.8333
$ws.Range("K107").Value = 972.8333
$ws.Range("M107").Value = 947.1667
$ws.Range("H134").Value = 4602.8
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -3465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23684.174
$ws.Range("I31").Value = 3318.05
$ws.Range("J31").Value = 91571.25
$ws.Range("K31").Value = 3318.05
$ws.Range("L31").Value = 91571.25
$ws.Range("M31").Value = -3023.05
$ws.Range("N31").Value = -92161.25
$ws.Range("H34").Value = 23684.174
$ws.Range("I34").Value = 3318.05
$ws.Range("J34").Value = 91571.25
$ws.Range("K34").Value = 3318.05
$ws.Range("L34").Value = 91571.25
$ws.Range("M34").Value = -3116.05
$ws.Range("N34").Value = -91975.25
$ws.Range("H86").Value = 17127
$ws.Range("I86").Value = 9500
$ws.Range("J86").Value = 19669.334
$ws.Range("K86").Value = 9500
$ws.Range("L86").Value = 19669.334
$ws.Range("M86").Value = -8377
$ws.Range("N86").Value = -21915.334
$ws.Range("H89").Value = 17127
$ws.Range("I89").Value = 9500
$ws.Range("J89").Value = 19669.334
$ws.Range("K89").Value = 47500
$ws.Range("L89").Value = 98346.67
$ws.Range("M89").Value = -41884
$ws.Range("N89").Value = -109578.67
$ws.Range("H122").Value = 6109.684
$ws.Range("I122").Value = 1775.7693
$ws.Range("J122").Value = 15499.833
$ws.Range("K122").Value = 5327.3079
$ws.Range("L122").Value = 46499.499
$ws.Range("M122").Value = -2877.3079
$ws.Range("N122").Value = -51399.499
$ws.Range("H134").Value = 3893.077
$ws.Range("I134").Value = 2474.75
$ws.Range("K134").Value = 7424.25
$ws.Range("M134").Value = -4889.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 7258.75
$ws.Range("J69").Value = 9506
$ws.Range("L69").Value = 28518
$ws.Range("N69").Value = -30140
$ws.Range("H70").Value = 14507
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 30000
$ws.Range("M70").Value = -29685
$ws.Range("H72").Value = 7258.75
$ws.Range("J72").Value = 9506
$ws.Range("L72").Value = 85554
$ws.Range("N72").Value = -93666
$ws.Range("H73").Value = 14507
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 30000
$ws.Range("M73").Value = -28908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1483.9166
$ws.Range("I97").Value = 970.4286
$ws.Range("K97").Value = 970.4286
$ws.Range("M97").Value = -474.4286
$ws.Range("H102").Value = 2684.625
$ws.Range("I102").Value = 1917.8334
$ws.Range("J102").Value = 4985
$ws.Range("K102").Value = 1917.8334
$ws.Range("L102").Value = 4985
$ws.Range("M102").Value = -295.8334
$ws.Range("N102").Value = -8229

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4463.273
$ws.Range("I46").Value = 1531.6666
$ws.Range("J46").Value = 5562.625
$ws.Range("K46").Value = 1531.6666
$ws.Range("L46").Value = 5562.625
$ws.Range("M46").Value = -1343.6666
$ws.Range("N46").Value = -5938.625
$ws.Range("H55").Value = 1896.125
$ws.Range("I55").Value = 487
$ws.Range("J55").Value = 3305.25
$ws.Range("K55").Value = 487
$ws.Range("L55").Value = 3305.25
$ws.Range("M55").Value = -314
$ws.Range("N55").Value = -3651.25
$ws.Range("H100").Value = 10986.546
$ws.Range("I100").Value = 10706.286
$ws.Range("K100").Value = 10706.286
$ws.Range("M100").Value = -10165.286
$ws.Range("H122").Value = 7952.25
$ws.Range("I122").Value = 3902
$ws.Range("K122").Value = 11706
$ws.Range("M122").Value = -9256
$ws.Range("H132").Value = 19005
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 19005
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 57015
$ws.Range("N132").Value = -62075
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 13490
$ws.Range("I136").Value = 4624.6665
$ws.Range("K136").Value = 13873.9995
$ws.Range("M136").Value = -11323.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8504.4
$ws.Range("I132").Value = 7820.84
$ws.Range("K132").Value = 23462.52
$ws.Range("M132").Value = -20932.52
$ws.Range("H136").Value = 6286.8
$ws.Range("I136").Value = 3607.875
$ws.Range("K136").Value = 10823.625
$ws.Range("M136").Value = -8273.625
